# "added solidworks and LRA motor reference"
#
# The underlying data change is to the "Resistance" row (row 4) of the
# components spec sheet: the existing tolerance values for B4/D4 get an
# explicit unit ("Ohm") appended, and a new third value is added in G4
# (lining up with the " G0832012 LRA" motor column that already has data
# in rows 1, 2, 3 and 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two existing resistance values to include units.
$ws.Range("B4").Value = "32 +- 20% Ohm"
$ws.Range("D4").Value = "6 +- 5% Ohm"

# Add the new resistance value for the LRA motor column.
$ws.Range("G4").Value = "22 +- 18% Ohm"

# Reflect the updated view state: the sheet is scrolled down a bit and the
# active selection moves from L5 back to A5.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("A5").Select()
